# Updated cryptos list on Tue Nov 19 18:40:31 UTC 2024 with GitHub Actions
#
# Note: several "Price" values in column D look like plain numbers
# (e.g. "5.50", "243.53") but must stay as literal text so formatting
# such as trailing/leading zeros is preserved exactly as scraped. They
# are written with a leading apostrophe, which is Excel's standard way
# of forcing a value to be stored as text instead of being parsed as a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.669.48'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '3.135.41'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''243.53'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '''618.12'
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("D8").Value = '''0.417'
$ws.Range("E8").Value = '  +12.29%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '3.131.76'
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("D11").Value = '''0.735'
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '''0.0000260'
$ws.Range("E13").Value = '  +4.89%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''34.80'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '92.969.40'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").Value = '''5.50'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '3.721.36'
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").Value = '3.124.77'
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("D19").Value = '''3.81'
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").Value = '''14.85'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("B21").Value = 'PEPE'
$ws.Range("C21").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D21").Value = '''0.0000210'
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = '''5.82'
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '''452.06'
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").Value = '''9.43'
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("D25").Value = '''5.86'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").Value = '''87.76'
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").Value = '''11.90'
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("D28").Value = '3.293.38'
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +5.73%  '
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '''0.170'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.226'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").Value = '''8.11'
$ws.Range("E35").Value = '  +3.50%  '
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").Value = '''26.28'
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").Value = '''4.02'
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").Value = '''485.06'
$ws.Range("E40").Value = '  -5.06%  '
$ws.Range("D41").Value = '''1.31'
$ws.Range("E41").Value = '  -3.77%  '
$ws.Range("D42").Value = '''3.53'
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("D43").Value = '''0.437'
$ws.Range("E43").Value = '  -4.27%  '
$ws.Range("D44").Value = '''23.10'
$ws.Range("E44").Value = '  +4.31%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '''161.59'
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("D47").Value = '''1.95'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("D48").Value = '''0.697'
$ws.Range("E48").Value = '  -5.02%  '
$ws.Range("D49").Value = '''1.40'
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").Value = '''0.0333'
$ws.Range("E50").Value = '  +3.54%  '
$ws.Range("E51").Value = '  +0.67%  '
